$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67634766666667
$ws.Range("H2").Value = 83.029043
$ws.Range("I2").Value = 0.005965811625935536
$ws.Range("J2").Value = 0.005965811625935536
$ws.Range("M2").Value = 164.855011
$ws.Range("N2").Value = 494.565033
$ws.Range("O2").Value = 0.8897130830256036
$ws.Range("P2").Value = 0.8897130830256037
$ws.Range("Q2").Value = 4562.584599028158
$ws.Range("R2").Value = 41063.26139125341
$ws.Range("S2").Value = 0.005307860654461095
$ws.Range("T2").Value = 0.005307860654461095
$ws.Range("G3").Value = 27.67634766666667
$ws.Range("H3").Value = 83.029043
$ws.Range("I3").Value = 0.005965811625935536
$ws.Range("J3").Value = 0.005965811625935536
$ws.Range("O3").Value = 0.009319291703989
$ws.Range("P3").Value = 0.009319291703989001
$ws.Range("Q3").Value = 47.79075143851455
$ws.Range("R3").Value = 430.116762946631
$ws.Range("S3").Value = 0.00005559713879314216
$ws.Range("T3").Value = 0.00005559713879314217
$ws.Range("G4").Value = 27.67634766666667
$ws.Range("H4").Value = 83.029043
$ws.Range("I4").Value = 0.005965811625935536
$ws.Range("J4").Value = 0.005965811625935536
$ws.Range("M4").Value = 4.66298
$ws.Range("N4").Value = 13.98894
$ws.Range("O4").Value = 0.02516583685701085
$ws.Range("P4").Value = 0.02516583685701086
$ws.Range("Q4").Value = 129.0542556427133
$ws.Range("R4").Value = 1161.48830078442
$ws.Range("S4").Value = 0.0001501346420979523
$ws.Range("T4").Value = 0.0001501346420979524
$ws.Range("G5").Value = 27.67634766666667
$ws.Range("H5").Value = 83.029043
$ws.Range("I5").Value = 0.005965811625935536
$ws.Range("J5").Value = 0.005965811625935536
$ws.Range("M5").Value = 14.04531966666667
$ws.Range("N5").Value = 42.135959
$ws.Range("O5").Value = 0.07580178841339644
$ws.Range("P5").Value = 0.07580178841339645
$ws.Range("Q5").Value = 388.7231501841375
$ws.Range("R5").Value = 3498.508351657237
$ws.Range("S5").Value = 0.000452219190583346
$ws.Range("T5").Value = 0.0004522191905833462
$ws.Range("I6").Value = 0.009118181457976757
$ws.Range("J6").Value = 0.009118181457976757
$ws.Range("M6").Value = 164.855011
$ws.Range("N6").Value = 494.565033
$ws.Range("O6").Value = 0.8897130830256036
$ws.Range("P6").Value = 0.8897130830256037
$ws.Range("Q6").Value = 6973.481044967612
$ws.Range("R6").Value = 62761.3294047085
$ws.Range("S6").Value = 0.008112565336563394
$ws.Range("T6").Value = 0.008112565336563394
$ws.Range("I7").Value = 0.009118181457976757
$ws.Range("J7").Value = 0.009118181457976757
$ws.Range("O7").Value = 0.009319291703989
$ws.Range("P7").Value = 0.009319291703989001
$ws.Range("S7").Value = 0.00008497499281678911
$ws.Range("T7").Value = 0.00008497499281678912
$ws.Range("I8").Value = 0.009118181457976757
$ws.Range("J8").Value = 0.009118181457976757
$ws.Range("M8").Value = 4.66298
$ws.Range("N8").Value = 13.98894
$ws.Range("O8").Value = 0.02516583685701085
$ws.Range("P8").Value = 0.02516583685701086
$ws.Range("Q8").Value = 197.24728078216
$ws.Range("R8").Value = 1775.22552703944
$ws.Range("S8").Value = 0.0002294666670040644
$ws.Range("T8").Value = 0.0002294666670040644
$ws.Range("I9").Value = 0.009118181457976757
$ws.Range("J9").Value = 0.009118181457976757
$ws.Range("M9").Value = 14.04531966666667
$ws.Range("N9").Value = 42.135959
$ws.Range("O9").Value = 0.07580178841339644
$ws.Range("P9").Value = 0.07580178841339645
$ws.Range("Q9").Value = 594.1267412612093
$ws.Range("R9").Value = 5347.140671350884
$ws.Range("S9").Value = 0.0006911744615925087
$ws.Range("T9").Value = 0.0006911744615925089
$ws.Range("G10").Value = 29.593002
$ws.Range("H10").Value = 88.779006
$ws.Range("I10").Value = 0.006378958578792732
$ws.Range("J10").Value = 0.006378958578792732
$ws.Range("M10").Value = 164.855011
$ws.Range("N10").Value = 494.565033
$ws.Range("O10").Value = 0.8897130830256036
$ws.Range("P10").Value = 0.8897130830256037
$ws.Range("Q10").Value = 4878.554670233021
$ws.Range("R10").Value = 43906.99203209719
$ws.Range("S10").Value = 0.005675442903630304
$ws.Range("T10").Value = 0.005675442903630305
$ws.Range("G11").Value = 29.593002
$ws.Range("H11").Value = 88.779006
$ws.Range("I11").Value = 0.006378958578792732
$ws.Range("J11").Value = 0.006378958578792732
$ws.Range("O11").Value = 0.009319291703989
$ws.Range("P11").Value = 0.009319291703989001
$ws.Range("Q11").Value = 51.10037711387799
$ws.Range("R11").Value = 459.903394024902
$ws.Range("S11").Value = 0.00005944737576343256
$ws.Range("T11").Value = 0.00005944737576343258
$ws.Range("G12").Value = 29.593002
$ws.Range("H12").Value = 88.779006
$ws.Range("I12").Value = 0.006378958578792732
$ws.Range("J12").Value = 0.006378958578792732
$ws.Range("M12").Value = 4.66298
$ws.Range("N12").Value = 13.98894
$ws.Range("O12").Value = 0.02516583685701085
$ws.Range("P12").Value = 0.02516583685701086
$ws.Range("Q12").Value = 137.99157646596
$ws.Range("R12").Value = 1241.92418819364
$ws.Range("S12").Value = 0.0001605318309115277
$ws.Range("T12").Value = 0.0001605318309115277
$ws.Range("G13").Value = 29.593002
$ws.Range("H13").Value = 88.779006
$ws.Range("I13").Value = 0.006378958578792732
$ws.Range("J13").Value = 0.006378958578792732
$ws.Range("M13").Value = 14.04531966666667
$ws.Range("N13").Value = 42.135959
$ws.Range("O13").Value = 0.07580178841339644
$ws.Range("P13").Value = 0.07580178841339645
$ws.Range("Q13").Value = 415.643172986306
$ws.Range("R13").Value = 3740.788556876754
$ws.Range("S13").Value = 0.0004835364684874667
$ws.Range("T13").Value = 0.0004835364684874668
$ws.Range("G14").Value = 4539.588785666667
$ws.Range("H14").Value = 13618.766357
$ws.Range("I14").Value = 0.9785370483372949
$ws.Range("J14").Value = 0.978537048337295
$ws.Range("M14").Value = 164.855011
$ws.Range("N14").Value = 494.565033
$ws.Range("O14").Value = 0.8897130830256036
$ws.Range("P14").Value = 0.8897130830256037
$ws.Range("Q14").Value = 748373.959196555
$ws.Range("R14").Value = 6735365.632768994
$ws.Range("S14").Value = 0.8706172141309487
$ws.Range("T14").Value = 0.870617214130949
$ws.Range("G15").Value = 4539.588785666667
$ws.Range("H15").Value = 13618.766357
$ws.Range("I15").Value = 0.9785370483372949
$ws.Range("J15").Value = 0.978537048337295
$ws.Range("O15").Value = 0.009319291703989
$ws.Range("P15").Value = 0.009319291703989001
$ws.Range("Q15").Value = 7838.836319799462
$ws.Range("R15").Value = 70549.52687819516
$ws.Range("S15").Value = 0.009119272196615635
$ws.Range("T15").Value = 0.009119272196615638
$ws.Range("G16").Value = 4539.588785666667
$ws.Range("H16").Value = 13618.766357
$ws.Range("I16").Value = 0.9785370483372949
$ws.Range("J16").Value = 0.978537048337295
$ws.Range("M16").Value = 4.66298
$ws.Range("N16").Value = 13.98894
$ws.Range("O16").Value = 0.02516583685701085
$ws.Range("P16").Value = 0.02516583685701086
$ws.Range("Q16").Value = 21168.01171578795
$ws.Range("R16").Value = 190512.1054420916
$ws.Range("S16").Value = 0.02462570371699731
$ws.Range("T16").Value = 0.02462570371699731
$ws.Range("G17").Value = 4539.588785666667
$ws.Range("H17").Value = 13618.766357
$ws.Range("I17").Value = 0.9785370483372949
$ws.Range("J17").Value = 0.978537048337295
$ws.Range("M17").Value = 14.04531966666667
$ws.Range("N17").Value = 42.135959
$ws.Range("O17").Value = 0.07580178841339644
$ws.Range("P17").Value = 0.07580178841339645
$ws.Range("Q17").Value = 63759.97564990348
$ws.Range("R17").Value = 573839.7808491314
$ws.Range("S17").Value = 0.07417485829273311
$ws.Range("T17").Value = 0.07417485829273314
